$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C19").Value = 26
$ws.Range("D19").Value = "Data Cleaning"

$ws.Range("C20").Value = 27
$ws.Range("D20").Value = "Excel Project"

$ws.Range("H20").Select()
